# Apply cached-value updates to Anima_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Source data was refreshed by a scheduled market-data runner; only the cached
# numeric values in columns H-N changed (no formulas are present in these sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23313.857
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 23313.857
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H38").Value = 371.5
$ws.Range("I38").Value = 265.8
$ws.Range("K38").Value = 797.4000000000001
$ws.Range("M38").Value = -425.4000000000001
$ws.Range("H58").Value = 1084.2858
$ws.Range("I58").Value = 931.6667
$ws.Range("K58").Value = 2795.0001
$ws.Range("M58").Value = -2645.0001
$ws.Range("H138").Value = 2077.0505
$ws.Range("I138").Value = 1456.2609
$ws.Range("J138").Value = 2615.849
$ws.Range("K138").Value = 4368.7827
$ws.Range("L138").Value = 7847.547
$ws.Range("M138").Value = 771.2173000000003
$ws.Range("N138").Value = -18127.547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4495.7144
$ws.Range("I63").Value = 3967.5454
$ws.Range("K63").Value = 3967.5454
$ws.Range("M63").Value = -3281.5454
$ws.Range("H66").Value = 4495.7144
$ws.Range("I66").Value = 3967.5454
$ws.Range("K66").Value = 19837.727
$ws.Range("M66").Value = -16405.727
$ws.Range("H80").Value = 39212.5
$ws.Range("J80").Value = 39212.5
$ws.Range("L80").Value = 39212.5
$ws.Range("N80").Value = -41208.5
$ws.Range("H83").Value = 39212.5
$ws.Range("J83").Value = 39212.5
$ws.Range("L83").Value = 117637.5
$ws.Range("N83").Value = -127621.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 18258.4
$ws.Range("J35").Value = 18258.4
$ws.Range("L35").Value = 18258.4
$ws.Range("N35").Value = -18878.4
$ws.Range("H82").Value = 25686.572
$ws.Range("J82").Value = 39181.875
$ws.Range("L82").Value = 39181.875
$ws.Range("N82").Value = -39947.875
$ws.Range("H85").Value = 25686.572
$ws.Range("J85").Value = 39181.875
$ws.Range("L85").Value = 39181.875
$ws.Range("N85").Value = -41833.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 7499.625
$ws.Range("J41").Value = 8999.5
$ws.Range("L41").Value = 8999.5
$ws.Range("N41").Value = -9855.5
$ws.Range("H51").Value = 19999
$ws.Range("J51").Value = 19999
$ws.Range("L51").Value = 19999
$ws.Range("N51").Value = -21471
$ws.Range("H58").Value = 881.5738
$ws.Range("I58").Value = 662.2439000000001
$ws.Range("J58").Value = 1331.2
$ws.Range("K58").Value = 662.2439000000001
$ws.Range("L58").Value = 1331.2
$ws.Range("M58").Value = -459.2439000000001
$ws.Range("N58").Value = -1737.2
$ws.Range("H60").Value = 10504
$ws.Range("J60").Value = 10504
$ws.Range("L60").Value = 10504
$ws.Range("N60").Value = -11526
$ws.Range("H61").Value = 19999
$ws.Range("J61").Value = 19999
$ws.Range("L61").Value = 19999
$ws.Range("N61").Value = -20695
$ws.Range("H68").Value = 22626.3
$ws.Range("J68").Value = 22626.3
$ws.Range("L68").Value = 22626.3
$ws.Range("N68").Value = -24124.3
$ws.Range("H71").Value = 22626.3
$ws.Range("J71").Value = 22626.3
$ws.Range("L71").Value = 67878.89999999999
$ws.Range("N71").Value = -75366.89999999999
$ws.Range("H74").Value = 28877
$ws.Range("J74").Value = 28877
$ws.Range("L74").Value = 28877
$ws.Range("N74").Value = -30625
$ws.Range("H77").Value = 28877
$ws.Range("J77").Value = 28877
$ws.Range("L77").Value = 86631
$ws.Range("N77").Value = -95367
$ws.Range("H134").Value = 10203.333
$ws.Range("I134").Value = 11644
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 34932
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -32397
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 881.5738
$ws.Range("I136").Value = 662.2439000000001
$ws.Range("J136").Value = 1331.2
$ws.Range("K136").Value = 1986.7317
$ws.Range("L136").Value = 3993.6
$ws.Range("M136").Value = 563.2682999999997
$ws.Range("N136").Value = -9093.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7156.5713
$ws.Range("H46").Value = 3894.6843
$ws.Range("J46").Value = 3894.6843
$ws.Range("L46").Value = 3894.6843
$ws.Range("N46").Value = -4206.6843
$ws.Range("H57").Value = 18999.908
$ws.Range("J57").Value = 18999.908
$ws.Range("L57").Value = 18999.908
$ws.Range("N57").Value = -20639.908
$ws.Range("H70").Value = 5433.7744
$ws.Range("I70").Value = 5358.609
$ws.Range("K70").Value = 5358.609
$ws.Range("M70").Value = -5088.609
$ws.Range("H73").Value = 5433.7744
$ws.Range("I73").Value = 5358.609
$ws.Range("K73").Value = 5358.609
$ws.Range("M73").Value = -4422.609
$ws.Range("H80").Value = 127350000
$ws.Range("I80").Value = 169666670
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 169666670
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -169665672
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 127350000
$ws.Range("I83").Value = 169666670
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 848333350
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -848328358
$ws.Range("N83").Value = -2009984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2373.818
$ws.Range("I132").Value = 2445.9092
$ws.Range("K132").Value = 7337.7276
$ws.Range("M132").Value = -4807.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3353.7778
$ws.Range("I81").Value = 3769.818
$ws.Range("K81").Value = 7539.636
$ws.Range("M81").Value = -6478.636
$ws.Range("H84").Value = 3353.7778
$ws.Range("I84").Value = 3769.818
$ws.Range("K84").Value = 37698.18
$ws.Range("M84").Value = -32394.18
$ws.Range("H107").Value = 654.0857
$ws.Range("I107").Value = 664.2963
$ws.Range("J107").Value = 619.625
$ws.Range("K107").Value = 1992.8889
$ws.Range("L107").Value = 1858.875
$ws.Range("M107").Value = -72.88889999999992
$ws.Range("N107").Value = -5698.875
$ws.Range("H110").Value = 330000
$ws.Range("J110").Value = 330000
$ws.Range("L110").Value = 330000
$ws.Range("N110").Value = -338180
$ws.Range("H132").Value = 2381787.2
$ws.Range("I132").Value = 716.6981
$ws.Range("J132").Value = 9805125
$ws.Range("K132").Value = 2150.0943
$ws.Range("L132").Value = 29415375
$ws.Range("M132").Value = 379.9057000000003
$ws.Range("N132").Value = -29420435
